$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume/1h changes, and two coin swaps)
$ws.Range("D2").Value = "'" + '52.410.40'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = "'" + '2.914.66'
$ws.Range("E3").Value = '  +3.93%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'" + '353.62'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").Value = "'" + '112.91'
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("D7").Value = "'" + '0.561'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'" + '0.633'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'" + '40.15'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = "'" + '0.0869'
$ws.Range("E11").Value = '  +3.62%  '
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").Value = "'" + '19.91'
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").Value = "'" + '7.84'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = "'" + '3.369.81'
$ws.Range("E15").Value = '  +3.88%  '
$ws.Range("E16").Value = '  +5.84%  '
$ws.Range("D17").Value = "'" + '2.909.22'
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("D18").Value = "'" + '52.421.62'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("B19").Value = 'ImmutableX'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D19").Value = "'" + '3.32'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'" + '7.64'
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = "'" + '14.25'
$ws.Range("E21").Value = '  +4.33%  '
$ws.Range("D22").Value = "'" + '0.0₃0982'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").Value = "'" + '71.03'
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("D24").Value = "'" + '270.42'
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("D25").Value = "'" + '2.79'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").Value = "'" + '0.174'
$ws.Range("E26").Value = '  +7.55%  '
$ws.Range("D27").Value = "'" + '26.87'
$ws.Range("E27").Value = '  +2.72%  '
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = "'" + '10.70'
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("D30").Value = "'" + '6.70'
$ws.Range("E30").Value = '  +9.22%  '
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Value = "'" + '6.39'
$ws.Range("E31").Value = '  +13.19%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = "'" + '37.95'
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").Value = "'" + '0.0979'
$ws.Range("E34").Value = '  +10.72%  '
$ws.Range("D35").Value = "'" + '53.29'
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("D36").Value = "'" + '0.0453'
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("D37").Value = "'" + '0.998'
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").Value = "'" + '3.35'
$ws.Range("E38").Value = '  +6.09%  '
$ws.Range("D39").Value = "'" + '18.91'
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("E40").Value = '  +3.27%  '
$ws.Range("D41").Value = "'" + '2.86'
$ws.Range("E41").Value = '  +14.55%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").Value = "'" + '23.55'
$ws.Range("E43").Value = '  +7.24%  '
$ws.Range("D44").Value = "'" + '121.40'
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("D45").Value = "'" + '2.62'
$ws.Range("E45").Value = '  +8.03%  '
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").Value = "'" + '3.57'
$ws.Range("E47").Value = '  +4.61%  '
$ws.Range("D48").Value = "'" + '2.203.52'
$ws.Range("E48").Value = '  +4.19%  '
$ws.Range("E49").Value = '  +21.55%  '
$ws.Range("D50").Value = "'" + '0.0343'
$ws.Range("E50").Value = '  +11.47%  '
$ws.Range("D51").Value = "'" + '0.972'
$ws.Range("E51").Value = '  +2.12%  '
